# Applies the "Horarios actualizados Linea 141 - 195" update:
#  - refreshes the "Ultima actualizacion" timestamp and "Total filas" count
#    on all three sheets
#  - inserts / appends new scraped rows reflecting the new scrape pass
#    (Hora_Scrap = 04:51:28) on sheet LP1912, LP1912-215 and 6203-6173

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:51:28"
$ws1.Range("A3").Value = "Total filas: 21"

# A new row was scraped that lands between the existing row 9 and row 10,
# pushing every following row down by one.
$ws1.Rows.Item(10).Insert()

$ws1.Cells.Item(10, 1).Value = "04:51:28"
$ws1.Cells.Item(10, 2).Value = "05:13"
$ws1.Cells.Item(10, 3).Value = "14_ABASTO"
$ws1.Cells.Item(10, 4).Value = 22
$ws1.Cells.Item(10, 5).Value = "LP1912"

# Two brand-new rows appended at the bottom of the sheet (rows 25 & 26).
$ws1.Cells.Item(25, 1).Value = "04:51:28"
$ws1.Cells.Item(25, 2).Value = "06:44"
$ws1.Cells.Item(25, 3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(25, 4).Value = 113
$ws1.Cells.Item(25, 5).Value = "LP1912"

$ws1.Cells.Item(26, 1).Value = "04:51:28"
$ws1.Cells.Item(26, 2).Value = "06:46"
$ws1.Cells.Item(26, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(26, 4).Value = 115
$ws1.Cells.Item(26, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:51:28"
$ws2.Range("A3").Value = "Total filas: 6"

# One brand-new row appended at the bottom of the sheet (row 11).
$ws2.Cells.Item(11, 1).Value = "04:51:28"
$ws2.Cells.Item(11, 2).Value = "06:46"
$ws2.Cells.Item(11, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(11, 4).Value = 115
$ws2.Cells.Item(11, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:51:28"
$ws3.Range("A3").Value = "Total filas: 6"

# A new row was scraped that lands between the existing row 7 and row 8,
# pushing the following rows down by one.
$ws3.Rows.Item(8).Insert()

$ws3.Cells.Item(8, 1).Value = "04:51:28"
$ws3.Cells.Item(8, 2).Value = "06:09"
$ws3.Cells.Item(8, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(8, 4).Value = 78
$ws3.Cells.Item(8, 5).Value = "L6173"

# One brand-new row appended at the bottom of the sheet (row 11).
$ws3.Cells.Item(11, 1).Value = "04:51:28"
$ws3.Cells.Item(11, 2).Value = "06:33"
$ws3.Cells.Item(11, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(11, 4).Value = 102
$ws3.Cells.Item(11, 5).Value = "L6203"
